$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column headers to remove spaces (replace with underscores)
$ws.Range("C1").Value = "Pre_Tell"
$ws.Range("D1").Value = "Gen_Pref"
$ws.Range("G1").Value = "Act_Pref"
